$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Guru99 Bank Manager HomePage"
$ws.Range("A3").HorizontalAlignment = -4108  # xlCenter

$ws.Range("A3").Select()
